$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lookup_table")

$ws.Range('B2').Value = '2033-2052'
$ws.Range('B3').Value = '2026-2045'
$ws.Range('B4').Value = '2022-2041'
$ws.Range('B5').Value = '2033-2052'
$ws.Range('B6').Value = '2007-2026'
$ws.Range('C6').Value = '2018-2037'
$ws.Range('B7').Value = '2016-2035'
$ws.Range('B8').Value = '2014-2033'
$ws.Range('C8').Value = '2028-2047'
$ws.Range('B9').Value = '2075-2094'
$ws.Range('B10').Value = '2039-2058'
$ws.Range('C10').Value = '2052-2071'
$ws.Range('B11').Value = '2034-2053'
$ws.Range('C11').Value = '2046-2065'
$ws.Range('B12').Value = '2025-2044'
$ws.Range('C12').Value = '2044-2063'
$ws.Range('B13').Value = '2032-2051'
$ws.Range('C13').Value = '2047-2066'
$ws.Range('B14').Value = '2019-2038'
$ws.Range('C14').Value = '2036-2055'
$ws.Range('B15').Value = '2015-2034'
$ws.Range('C15').Value = '2054-2073'
$ws.Range('B16').Value = '2023-2042'
$ws.Range('C16').Value = '2056-2075'
$ws.Range('B17').Value = '2024-2043'
$ws.Range('C17').Value = '2052-2071'
$ws.Range('B18').Value = '2022-2041'
$ws.Range('C18').Value = '2034-2053'
$ws.Range('D18').Value = '2069-2088'
$ws.Range('B19').Value = '2010-2029'
$ws.Range('C19').Value = '2024-2043'
$ws.Range('D19').Value = '2051-2070'
$ws.Range('B20').Value = '2007-2026'
$ws.Range('C20').Value = '2028-2047'
$ws.Range('D20').Value = '2054-2073'
$ws.Range('B21').Value = '2023-2042'
$ws.Range('C21').Value = '2071-2090'
$ws.Range('B22').Value = '2019-2038'
$ws.Range('C22').Value = '2029-2048'
$ws.Range('B23').Value = '2026-2045'
$ws.Range('C23').Value = '2043-2062'
$ws.Range('D23').Value = '2062-2081'
$ws.Range('B24').Value = '2019-2038'
$ws.Range('C24').Value = '2035-2054'
$ws.Range('D24').Value = '2055-2074'
$ws.Range('E24').Value = '2069-2088'
$ws.Range('B25').Value = '2019-2038'
$ws.Range('C25').Value = '2031-2050'
$ws.Range('D25').Value = '2049-2068'
$ws.Range('E25').Value = '2066-2085'
$ws.Range('B26').Value = '2024-2043'
$ws.Range('C26').Value = '2037-2056'
$ws.Range('D26').Value = '2065-2084'
$ws.Range('E26').Value = '2076-2095'
$ws.Range('B27').Value = '2015-2034'
$ws.Range('C27').Value = '2024-2043'
$ws.Range('D27').Value = '2048-2067'
$ws.Range('E27').Value = '2065-2084'
$ws.Range('B28').Value = '2015-2034'
$ws.Range('C28').Value = '2046-2065'
$ws.Range('D28').Value = '2062-2081'
$ws.Range('E28').Value = '2072-2091'
$ws.Range('B29').Value = '2015-2034'
$ws.Range('C29').Value = '2033-2052'
$ws.Range('D29').Value = '2065-2084'
$ws.Range('E29').Value = '2078-2097'
$ws.Range('B30').Value = '2038-2057'
$ws.Range('C30').Value = '2044-2063'
$ws.Range('D30').Value = '2060-2079'
$ws.Range('B31').Value = '2015-2034'
$ws.Range('C31').Value = '2031-2050'
$ws.Range('D31').Value = '2045-2064'
$ws.Range('E31').Value = '2057-2076'
$ws.Range('B32').Value = '2006-2025'
$ws.Range('C32').Value = '2019-2038'
$ws.Range('D32').Value = '2043-2062'
$ws.Range('E32').Value = '2055-2074'
$ws.Range('B33').Value = '2006-2025'
$ws.Range('C33').Value = '2022-2041'
$ws.Range('D33').Value = '2045-2064'
$ws.Range('E33').Value = '2056-2075'
$ws.Range('B34').Value = '2025-2044'
$ws.Range('C34').Value = '2040-2059'
$ws.Range('D34').Value = '2066-2085'
$ws.Range('E34').Value = '2080-2099'
$ws.Range('B35').Value = '2025-2044'
$ws.Range('C35').Value = '2036-2055'
$ws.Range('D35').Value = '2053-2072'
$ws.Range('E35').Value = '2067-2086'

$ws.Range("I34").Select()
